# Malingas workbook update:
#  - Insert a new row 12 "PIÑIN MACHUCA GIANCARLO" in both worksheets
#    (crosstab + annot), pushing the existing rows 12-15 down to 13-16.
#  - Update several "04dec2025" (column E) values across the table.
#
# Sheet "crosstab" stores genuine numbers; sheet "annot" mirrors the same
# data but as text, with numeric 0 shown as a blank cell.

$wb = $excel.ActiveWorkbook

$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J")

# ---------------------------------------------------------------------
# Final data for every data row (2-16) of the table, column order
# B..J = 01dec2025, 02dec2025, 03dec2025, 04dec2025, 26nov2025,
#        27nov2025, 28nov2025, 29nov2025, 30nov2025
# ---------------------------------------------------------------------
$names = @{
    2  = "AGURTO ORDINOLA LISBET JAQUELIN"
    3  = "ALAMA NIMA CLARITZA MABEL"
    4  = "ALBIRENA GARCIA ANGEELO ALONSO"
    5  = "ATOCHE PALACIOS LUIS ANGEL"
    6  = "BERNAOLA CARMEN ZUMIKO YASHURY"
    7  = "CARREÑO PALACIOS KATHERINE DE LOS MILAGROS"
    8  = "CORDOVA CARMEN ANGIE NATALLY"
    9  = "JUAREZ CARMEN PIERRE ALEXANDER"
    10 = "MANUEL LEUNARDO PRADO BAILON"
    11 = "MARYURI OJEDA VALLE"
    12 = "PIÑIN MACHUCA GIANCARLO"
    13 = "ROMAN GALECIO MARITZA DEL PILAR"
    14 = "RUIDIAS FRIAS MELISSA VICTORIA"
    15 = "URRIOLA ARISMENDIZ INGRID MARYURI"
    16 = "VEGA ROBLEDO FERNANDO ERNESTO"
}

$values = @{
    2  = @(33, 17, 0, 55, 14, 15, 27, 13, 0)
    3  = @(15, 15, 21, 22, 12, 17, 27, 20, 30)
    4  = @(21, 11, 3, 48, 14, 15, 35, 13, 48)
    5  = @(18, 33, 30, 27, 0, 15, 19, 24, 16)
    6  = @(15, 18, 25, 24, 0, 15, 15, 23, 8)
    7  = @(22, 8, 0, 20, 15, 0, 41, 0, 14)
    8  = @(6, 28, 22, 28, 10, 17, 15, 18, 40)
    9  = @(39, 17, 24, 22, 0, 25, 15, 15, 0)
    10 = @(9, 14, 0, 46, 14, 15, 19, 16, 60)
    11 = @(16, 31, 33, 24, 13, 15, 0, 40, 12)
    12 = @(0, 0, 0, 1, 0, 0, 0, 0, 0)
    13 = @(17, 15, 34, 16, 11, 16, 24, 17, 0)
    14 = @(9, 31, 0, 65, 10, 15, 15, 24, 0)
    15 = @(6, 26, 26, 35, 13, 18, 16, 27, 30)
    16 = @(8, 15, 25, 23, 0, 25, 23, 13, 27)
}

for ($s = 1; $s -le 2; $s++) {
    $ws = $wb.Worksheets.Item($s)

    # Insert the new row 12, pushing old rows 12-15 down to 13-16, and
    # clone the formatting of the row directly above it so the new row
    # keeps the same borders / bold header-style first column.
    $ws.Rows.Item(12).Insert()
    $ws.Range("A11:J11").Copy()
    $ws.Range("A12:J12").PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    if ($s -eq 2) {
        # "annot" mirrors the numbers as literal text, with 0 shown blank.
        $ws.Range("B2:J16").NumberFormat = "@"
    }

    foreach ($r in 2..16) {
        $ws.Range("A" + $r).Value = $names[$r]
        $rowVals = $values[$r]
        for ($i = 0; $i -lt 9; $i++) {
            $cellRef = $cols[$i] + $r
            $n = $rowVals[$i]
            if ($s -eq 1) {
                $ws.Range($cellRef).Value = $n
            } else {
                if ($n -eq 0) {
                    $ws.Range($cellRef).Value = ""
                } else {
                    $ws.Range($cellRef).Value = [string]$n
                }
            }
        }
    }
}
